$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# New shared strings must be introduced in this exact order so the saved
# sharedStrings table lines up with the target workbook (29..33):
#   1) "更改刷新逻辑" (A23)  2) "没有必要性" (E14)  3) "待定" (C21)
#   4) "找不到合适的界面设计" (E21)  5) "版本信息" (A24)

# --- Row 23 (new): change refresh logic --------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A23").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A23").Value = "更改刷新逻辑"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B23").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B23").Value = 42805

$ws.Rows.Item(23).RowHeight = 31.5

# --- Row 14: fill in status ("x") and remark ("not necessary") -------------
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C14").Value = "×"

$ws.Range("E17").Copy() | Out-Null
$ws.Range("E14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E14").Value = "没有必要性"

# --- Row 21: fill in status ("pending") and remark --------------------------
$ws.Range("C21").Value = "待定"

$ws.Range("E17").Copy() | Out-Null
$ws.Range("E21").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E21").Value = "找不到合适的界面设计"

# --- Row 24 (new): version info ----------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A24").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A24").Value = "版本信息"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B24").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B24").Value = 42805

$ws.Rows.Item(24).RowHeight = 33.75

# --- Row 20: fill in resolved date -----------------------------------------
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D20").Value = 42793

# --- Row 22: fill in proposed/resolved dates and status ---------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B22").Value = 42803

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C22").Value = "√"

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D22").Value = 42804

$excel.CutCopyMode = $false

# --- Update view: scroll down and select D24 --------------------------------
$ws.Range("D24").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
